$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-08"

# Update the header label in I1 to reflect the new "through" date
$ws.Range("I1").Value = "2022 (through 10-08)"

# Update the monthly total column (I) with revised figures
$ws.Range("I9").Value = 164
$ws.Range("I10").Value = 146
$ws.Range("I11").Value = 28
$ws.Range("I14").Value = 1306
